# Update "想去人数" (F column) figures across the relevant worksheets.
# Values were refreshed as of the output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7065
$ws1.Range("F4").Value  = 77
$ws1.Range("F7").Value  = 7602
$ws1.Range("F13").Value = 432
$ws1.Range("F14").Value = 160
$ws1.Range("F16").Value = 425
$ws1.Range("F20").Value = 5452
$ws1.Range("F21").Value = 139
$ws1.Range("F22").Value = 195
$ws1.Range("F23").Value = 880
$ws1.Range("F24").Value = 224
$ws1.Range("F25").Value = 291

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7065
$ws4.Range("F4").Value  = 77
$ws4.Range("F7").Value  = 7602
$ws4.Range("F13").Value = 432
$ws4.Range("F14").Value = 160
$ws4.Range("F16").Value = 425
$ws4.Range("F20").Value = 2
$ws4.Range("F21").Value = 5452
$ws4.Range("F23").Value = 139
$ws4.Range("F24").Value = 195
$ws4.Range("F25").Value = 880
$ws4.Range("F26").Value = 224
$ws4.Range("F27").Value = 291
